# kassir report excel file changed
#
# The "TDSheet" report gains a new "Тўлов тури" (payment type) column,
# inserted right before the existing "Тўланган сумма" (amount paid) column.
# Everything from the old column K ("Тўланган сумма") through column Q
# ("Реферал тел раками") shifts one column to the right (L..R), the hidden
# AutoFilter range grows to match, and the view is left scrolled over to
# show the new column with the cursor sitting on the first data row below
# the new header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at K - this pushes the old K:Q ("Тўланган сумма"
# .. "Реферал тел раками") columns one slot to the right, to L:R, carrying
# their widths/styles/values along with them.
$ws.Columns("K").Insert()

# Give the freshly inserted column K its header text on the header row.
$ws.Range("K4").Value = "Тўлов тури"

# The worksheet's hidden AutoFilter defined name covered A4:Q4; extend it
# by one column so it still spans the whole header row (now A4:R4).
$wb.Names.Item(1).RefersTo = "=TDSheet!`$A`$4:`$R`$4"

# Match the saved view: scrolled so column F is left-most on screen, with
# the active cell resting on K5 (just under the new header).
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K5").Select()
